# Corrected code for g expl & the arguments for the code in the separate excel file.
# Applies Tim's feedback: re-codes several "g" (geographic-limitation) explanations
# in column M ("Coding explanation"), flips the corresponding expl_c / expl_g
# boolean flags (columns F / J) where the underlying reason category changed,
# and refreshes the sheet view (freeze header row, cursor parked at M86).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param(
        [int]$Row,
        [string]$M,
        $F = $null,
        $J = $null
    )
    $ws.Range("M$Row").Value = $M
    if ($null -ne $F) { $ws.Range("F$Row").Value = $F }
    if ($null -ne $J) { $ws.Range("J$Row").Value = $J }
}

# positional args: Row, M, F, J  (F/J omitted => left unchanged)
Set-Row 18  "licensing lags e; MRV - 2015" 1 0
Set-Row 48  "Indicator is updated using raw data from 2+ years prior - MRV 2017 c;" $null 0
Set-Row 50  "MRV 2018" $null 0
Set-Row 53  "MRV 2018" $null 0
Set-Row 54  "MRV 2018" $null 0
Set-Row 55  "MRV 2018" $null 0
Set-Row 56  "MRV 2018" $null 0
Set-Row 57  "MRV 2018" $null 0
Set-Row 84  "MRV 2015 - c" 1 0
Set-Row 85  "MRV 2015 - c" 1 0
Set-Row 104 "MRV 2017 - c" $null 0
Set-Row 111 "MRV 2018" $null 0
Set-Row 116 "MRV 2018" $null 0
Set-Row 129 "MRV 2018" $null 0

# Refresh the sheet view: freeze the header row and leave the cursor on M86,
# matching the saved view state in the corrected workbook.
$ws.Activate()
$ws.Range("D1").Select()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("M86").Select()
